$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the current row 213 (pushes old row 214 -> 216)
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# Row 213: update with the new weekly reading (Primera) - date + prices change
$ws.Cells.Item(213,4).Value = 44628
$ws.Cells.Item(213,10).Value = 220
$ws.Cells.Item(213,11).Value = 9000
$ws.Cells.Item(213,12).Value = 10000
$ws.Cells.Item(213,13).Value = 9545
$ws.Cells.Item(213,16).Value = 1591

# Row 214 (new): new weekly reading (Segunda) for the same date
$ws.Cells.Item(214,1).Value = 11
$ws.Cells.Item(214,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(214,3).Value = "Bíobío"
$ws.Cells.Item(214,4).Value = 44628
$ws.Cells.Item(214,5).Value = 8
$ws.Cells.Item(214,6).Value = 100112017
$ws.Cells.Item(214,7).Value = "Apio"
$ws.Cells.Item(214,8).Value = "Americana (o)"
$ws.Cells.Item(214,9).Value = "Segunda"
$ws.Cells.Item(214,10).Value = 220
$ws.Cells.Item(214,11).Value = 7000
$ws.Cells.Item(214,12).Value = 8000
$ws.Cells.Item(214,13).Value = 7545
$ws.Cells.Item(214,14).Value = "`$/docena de matas"
$ws.Cells.Item(214,15).Value = "Región de Coquimbo"
$ws.Cells.Item(214,16).Value = 1258
$ws.Cells.Item(214,17).Value = 6
$ws.Cells.Item(214,18).Value = "Hortaliza"

# Row 215 (new): re-insert the prior (now historical) Primera reading that used to live on row 213
$ws.Cells.Item(215,1).Value = 11
$ws.Cells.Item(215,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(215,3).Value = "Bíobío"
$ws.Cells.Item(215,4).Value = 44552
$ws.Cells.Item(215,5).Value = 8
$ws.Cells.Item(215,6).Value = 100112017
$ws.Cells.Item(215,7).Value = "Apio"
$ws.Cells.Item(215,8).Value = "Americana (o)"
$ws.Cells.Item(215,9).Value = "Primera"
$ws.Cells.Item(215,10).Value = 100
$ws.Cells.Item(215,11).Value = 7000
$ws.Cells.Item(215,12).Value = 8000
$ws.Cells.Item(215,13).Value = 7500
$ws.Cells.Item(215,14).Value = "`$/docena de matas"
$ws.Cells.Item(215,15).Value = "Región de Coquimbo"
$ws.Cells.Item(215,16).Value = 1250
$ws.Cells.Item(215,17).Value = 6
$ws.Cells.Item(215,18).Value = "Hortaliza"

# Row 216 is the old row 214 (Segunda, 44552) shifted down by the inserts above - values unchanged.
